$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Presure Sensors")
$ws.Range("A1").Value = "test"
